# Auto-generated edit script: applies per-cell numeric updates to the
# "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H:N) that
# the scheduled market-data runner refreshed for each crafting sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H19").Value = 1645557.5
$ws.Range("I19").Value = 2632175
$ws.Range("J19").Value = 1195
$ws.Range("K19").Value = 2632175
$ws.Range("L19").Value = 1195
$ws.Range("M19").Value = -2632000
$ws.Range("N19").Value = -1545
$ws.Range("H33").Value = 348.42856
$ws.Range("I33").Value = 70
$ws.Range("J33").Value = 459.8
$ws.Range("K33").Value = 70
$ws.Range("L33").Value = 459.8
$ws.Range("M33").Value = 159
$ws.Range("N33").Value = -917.8
$ws.Range("H44").Value = 15831.111
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 15831.111
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 15831.111
$ws.Range("N44").Value = -16755.111
$ws.Range("H64").Value = 3200
$ws.Range("I64").Value = 3000
$ws.Range("H67").Value = 3200
$ws.Range("I67").Value = 3000
$ws.Range("H70").Value = 2946.8696
$ws.Range("I70").Value = 2001.3334
$ws.Range("J70").Value = 3088.7
$ws.Range("K70").Value = 6004.0002
$ws.Range("L70").Value = 9266.099999999999
$ws.Range("M70").Value = -5734.0002
$ws.Range("N70").Value = -9806.099999999999
$ws.Range("H73").Value = 2946.8696
$ws.Range("I73").Value = 2001.3334
$ws.Range("J73").Value = 3088.7
$ws.Range("K73").Value = 6004.0002
$ws.Range("L73").Value = 9266.099999999999
$ws.Range("M73").Value = -5068.0002
$ws.Range("N73").Value = -11138.1
$ws.Range("H86").Value = 1612.1111
$ws.Range("I86").Value = 1381.5
$ws.Range("J86").Value = 1796.6
$ws.Range("K86").Value = 1381.5
$ws.Range("L86").Value = 1796.6
$ws.Range("M86").Value = -258.5
$ws.Range("N86").Value = -4042.6
$ws.Range("H89").Value = 1612.1111
$ws.Range("I89").Value = 1381.5
$ws.Range("J89").Value = 1796.6
$ws.Range("K89").Value = 6907.5
$ws.Range("L89").Value = 8983
$ws.Range("M89").Value = -1291.5
$ws.Range("N89").Value = -20215
$ws.Range("H113").Value = 4195.9
$ws.Range("I113").Value = 2694.4
$ws.Range("J113").Value = 5697.4
$ws.Range("K113").Value = 2694.4
$ws.Range("L113").Value = 5697.4
$ws.Range("M113").Value = 559.5999999999999
$ws.Range("N113").Value = -12205.4
$ws.Range("H115").Value = 1730.8125
$ws.Range("I115").Value = 1069.8572
$ws.Range("J115").Value = 2244.889
$ws.Range("K115").Value = 3209.5716
$ws.Range("L115").Value = 6734.667
$ws.Range("M115").Value = -1642.5716
$ws.Range("N115").Value = -9868.667000000001
$ws.Range("H137").Value = 3146.6897
$ws.Range("I137").Value = 1791.375
$ws.Range("J137").Value = 4814.769
$ws.Range("K137").Value = 5374.125
$ws.Range("L137").Value = 14444.307
$ws.Range("M137").Value = -2824.125
$ws.Range("N137").Value = -19544.307
$ws.Range("H138").Value = 4871.643
$ws.Range("I138").Value = 813.7273
$ws.Range("J138").Value = 6931.8154
$ws.Range("K138").Value = 2441.1819
$ws.Range("L138").Value = 20795.4462
$ws.Range("M138").Value = 2698.8181
$ws.Range("N138").Value = -31075.4462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5547.5293
$ws.Range("I32").Value = 4428
$ws.Range("J32").Value = 12040.8
$ws.Range("K32").Value = 4428
$ws.Range("L32").Value = 12040.8
$ws.Range("M32").Value = -4141
$ws.Range("N32").Value = -12614.8
$ws.Range("H61").Value = 1013.625
$ws.Range("I61").Value = 747
$ws.Range("J61").Value = 2169
$ws.Range("K61").Value = 747
$ws.Range("L61").Value = 2169
$ws.Range("M61").Value = -535
$ws.Range("N61").Value = -2593
$ws.Range("H102").Value = 1833.8235
$ws.Range("I102").Value = 1850.6364
$ws.Range("J102").Value = 1803
$ws.Range("K102").Value = 1850.6364
$ws.Range("L102").Value = 1803
$ws.Range("M102").Value = -228.6364000000001
$ws.Range("N102").Value = -5047
$ws.Range("H122").Value = 2316.1765
$ws.Range("I122").Value = 1308.6364
$ws.Range("J122").Value = 4163.3335
$ws.Range("K122").Value = 3925.9092
$ws.Range("L122").Value = 12490.0005
$ws.Range("M122").Value = -1475.9092
$ws.Range("N122").Value = -17390.0005
$ws.Range("H132").Value = 3850.6428
$ws.Range("I132").Value = 1687.4286
$ws.Range("J132").Value = 6013.857
$ws.Range("K132").Value = 5062.2858
$ws.Range("L132").Value = 18041.571
$ws.Range("M132").Value = -2532.2858
$ws.Range("N132").Value = -23101.571
$ws.Range("H136").Value = 1013.625
$ws.Range("I136").Value = 747
$ws.Range("J136").Value = 2169
$ws.Range("K136").Value = 2241
$ws.Range("L136").Value = 6507
$ws.Range("M136").Value = 309
$ws.Range("N136").Value = -11607

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7169.5713
$ws.Range("I20").Value = 2567.4
$ws.Range("J20").Value = 12479.77
$ws.Range("K20").Value = 2567.4
$ws.Range("L20").Value = 12479.77
$ws.Range("M20").Value = -2320.4
$ws.Range("N20").Value = -12973.77
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -827
$ws.Range("H134").Value = 1971.8788
$ws.Range("I134").Value = 1298.963
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 3896.889
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -1361.889
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 587.04
$ws.Range("I22").Value = 501
$ws.Range("J22").Value = 590.625
$ws.Range("K22").Value = 501
$ws.Range("L22").Value = 590.625
$ws.Range("M22").Value = -151
$ws.Range("N22").Value = -1290.625
$ws.Range("H31").Value = 13160104
$ws.Range("I31").Value = 967.24
$ws.Range("J31").Value = 38466136
$ws.Range("K31").Value = 967.24
$ws.Range("L31").Value = 38466136
$ws.Range("M31").Value = -672.24
$ws.Range("N31").Value = -38466726
$ws.Range("H34").Value = 13160104
$ws.Range("I34").Value = 967.24
$ws.Range("J34").Value = 38466136
$ws.Range("K34").Value = 967.24
$ws.Range("L34").Value = 38466136
$ws.Range("M34").Value = -765.24
$ws.Range("N34").Value = -38466540
$ws.Range("H58").Value = 1613.0952
$ws.Range("I58").Value = 1516.7313
$ws.Range("J58").Value = 1992.8823
$ws.Range("K58").Value = 1516.7313
$ws.Range("L58").Value = 1992.8823
$ws.Range("M58").Value = -1313.7313
$ws.Range("N58").Value = -2398.8823
$ws.Range("H132").Value = 3653.7222
$ws.Range("I132").Value = 3041.48
$ws.Range("J132").Value = 5045.1816
$ws.Range("K132").Value = 9124.440000000001
$ws.Range("L132").Value = 15135.5448
$ws.Range("M132").Value = -6594.440000000001
$ws.Range("N132").Value = -20195.5448
$ws.Range("H134").Value = 13152
$ws.Range("I134").Value = 21102.4
$ws.Range("J134").Value = 5201.6
$ws.Range("K134").Value = 63307.2
$ws.Range("L134").Value = 15604.8
$ws.Range("M134").Value = -60772.2
$ws.Range("N134").Value = -20674.8
$ws.Range("H136").Value = 1613.0952
$ws.Range("I136").Value = 1516.7313
$ws.Range("J136").Value = 1992.8823
$ws.Range("K136").Value = 4550.1939
$ws.Range("L136").Value = 5978.6469
$ws.Range("M136").Value = -2000.1939
$ws.Range("N136").Value = -11078.6469

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3695
$ws.Range("I117").Value = 365
$ws.Range("J117").Value = 15350
$ws.Range("K117").Value = 1095
$ws.Range("L117").Value = 46050
$ws.Range("M117").Value = 2347
$ws.Range("N117").Value = -52934
$ws.Range("H131").Value = 15152404
$ws.Range("I131").Value = 100000290
$ws.Range("J131").Value = 996.6070999999999
$ws.Range("K131").Value = 300000870
$ws.Range("L131").Value = 2989.8213
$ws.Range("M131").Value = -299995830
$ws.Range("N131").Value = -13069.8213

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 41669068
$ws.Range("I80").Value = 50002284
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 50002284
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -50001286
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 41669068
$ws.Range("I83").Value = 50002284
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 250011420
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -250006428
$ws.Range("N83").Value = -24984
$ws.Range("H122").Value = 2539.1765
$ws.Range("I122").Value = 1743.6154
$ws.Range("J122").Value = 5124.75
$ws.Range("K122").Value = 5230.8462
$ws.Range("L122").Value = 15374.25
$ws.Range("M122").Value = -2780.8462
$ws.Range("N122").Value = -20274.25
$ws.Range("H132").Value = 2595.7932
$ws.Range("I132").Value = 1540.909
$ws.Range("J132").Value = 5911.143
$ws.Range("K132").Value = 4622.727000000001
$ws.Range("L132").Value = 17733.429
$ws.Range("M132").Value = -2092.727000000001
$ws.Range("N132").Value = -22793.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 42000
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 50000
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 50000
$ws.Range("M74").Value = -9002
$ws.Range("N74").Value = -51996
$ws.Range("H77").Value = 42000
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 50000
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 150000
$ws.Range("M77").Value = -25008
$ws.Range("N77").Value = -159984
$ws.Range("H136").Value = 2276.6943
$ws.Range("I136").Value = 1202.1786
$ws.Range("J136").Value = 6037.5
$ws.Range("K136").Value = 3606.5358
$ws.Range("L136").Value = 18112.5
$ws.Range("M136").Value = -1056.5358
$ws.Range("N136").Value = -23212.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1440
$ws.Range("I81").Value = 1600
$ws.Range("J81").Value = 1066.6666
$ws.Range("K81").Value = 3200
$ws.Range("L81").Value = 2133.3332
$ws.Range("M81").Value = -2139
$ws.Range("N81").Value = -4255.3332
$ws.Range("H84").Value = 1440
$ws.Range("I84").Value = 1600
$ws.Range("J84").Value = 1066.6666
$ws.Range("K84").Value = 16000
$ws.Range("L84").Value = 10666.666
$ws.Range("M84").Value = -10696
$ws.Range("N84").Value = -21274.666
$ws.Range("H113").Value = 415.4
$ws.Range("I113").Value = 297.42856
$ws.Range("J113").Value = 518.625
$ws.Range("K113").Value = 892.28568
$ws.Range("L113").Value = 1555.875
$ws.Range("M113").Value = 1277.71432
$ws.Range("N113").Value = -5895.875
$ws.Range("H132").Value = 7753635
$ws.Range("I132").Value = 956.70966
$ws.Range("J132").Value = 27781388
$ws.Range("K132").Value = 2870.12898
$ws.Range("L132").Value = 83344164
$ws.Range("M132").Value = -340.12898
$ws.Range("N132").Value = -83349224
$ws.Range("H136").Value = 3617.8262
$ws.Range("I136").Value = 917.2222
$ws.Range("J136").Value = 13340
$ws.Range("K136").Value = 2751.6666
$ws.Range("L136").Value = 40020
$ws.Range("M136").Value = -201.6666
$ws.Range("N136").Value = -45120
